$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns (Price / Volume) stay text-formatted so
# numeric-looking strings (e.g. "0.9986", "240.04") are not
# auto-converted to numbers by Excel, matching the source data
# which stores these as inline strings.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.887.02"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "1.808.13"
$ws.Range("E3").Value = "  -1.22%  "
$ws.Range("D4").Value = "0.9986"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "240.04"
$ws.Range("E5").Value = "  -1.40%  "
$ws.Range("D6").Value = "0.6061"
$ws.Range("E6").Value = "  -3.45%  "
$ws.Range("D7").Value = "0.9987"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "0.07257"
$ws.Range("E8").Value = "  -3.29%  "
$ws.Range("D9").Value = "0.2858"
$ws.Range("E9").Value = "  -2.03%  "
$ws.Range("D10").Value = "22.73"
$ws.Range("E10").Value = "  -1.65%  "
$ws.Range("D11").Value = "0.07617"
$ws.Range("E11").Value = "  -1.24%  "
$ws.Range("D12").Value = "1.791.44"
$ws.Range("E12").Value = "  -2.18%  "
$ws.Range("D13").Value = "4.904"
$ws.Range("E13").Value = "  -1.66%  "
$ws.Range("D14").Value = "0.6543"
$ws.Range("E14").Value = "  -1.88%  "
$ws.Range("D15").Value = "80.74"
$ws.Range("E15").Value = "  -2.16%  "
$ws.Range("D16").Value = "0.000008950"
$ws.Range("E16").Value = "  -4.35%  "
$ws.Range("E17").Value = "  -2.97%  "
$ws.Range("D18").Value = "28.868.67"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("D19").Value = "2.064.48"
$ws.Range("E19").Value = "  -0.85%  "
$ws.Range("D20").Value = "236.05"
$ws.Range("E20").Value = "  +5.65%  "
$ws.Range("E21").Value = "  -1.98%  "
$ws.Range("D22").Value = "0.9991"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").Value = "7.065"
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("D24").Value = "0.9993"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").Value = "157.62"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("D26").Value = "0.1394"
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("D27").Value = "8.354"
$ws.Range("E27").Value = "  -1.57%  "
$ws.Range("D28").Value = "17.48"
$ws.Range("E28").Value = "  -2.33%  "
$ws.Range("D29").Value = "1.471"
$ws.Range("E29").Value = "  -1.80%  "
$ws.Range("E30").Value = "  -2.48%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "4.046"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "4.056"
$ws.Range("E32").Value = "  -2.24%  "
$ws.Range("D33").Value = "1.202"
$ws.Range("E33").Value = "  +0.24%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "0.7291"
$ws.Range("E34").Value = "  -2.41%  "
$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "1.805"
$ws.Range("E35").Value = "  -1.88%  "
$ws.Range("D36").Value = "1.123"
$ws.Range("E36").Value = "  -0.93%  "
$ws.Range("D37").Value = "2.620"
$ws.Range("E37").Value = "  -1.59%  "
$ws.Range("D38").Value = "2.802"
$ws.Range("E38").Value = "  +1.50%  "
$ws.Range("D39").Value = "0.01741"
$ws.Range("E39").Value = "  -2.28%  "
$ws.Range("D40").Value = "1.185.69"
$ws.Range("D41").Value = "6.309"
$ws.Range("E41").Value = "  -3.50%  "
$ws.Range("D42").Value = "0.8833"
$ws.Range("E42").Value = "  -0.92%  "
$ws.Range("D43").Value = "0.9979"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "99.95"
$ws.Range("E44").Value = "  -2.16%  "
$ws.Range("D45").Value = "1.938.49"
$ws.Range("E45").Value = "  -2.04%  "
$ws.Range("D46").Value = "0.00000000122"
$ws.Range("E46").Value = "  +1.00%  "
$ws.Range("D47").Value = "63.78"
$ws.Range("E47").Value = "  -2.53%  "
$ws.Range("D48").Value = "0.5070"
$ws.Range("E48").Value = "  -0.38%  "
$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D49").Value = "0.3957"
$ws.Range("E49").Value = "  -2.65%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "8.967"
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("D51").Value = "0.05770"
$ws.Range("E51").Value = "  -1.01%  "
